$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are plain decimal numbers need to stay as text
# (matching the source inline-string cells), so force a Text number format
# before assigning, otherwise Excel auto-converts them to numeric values
# and trailing zeros (e.g. "2.00" -> 2) would be lost.
$textCells = @('D5', 'D6', 'D10', 'D12', 'D16', 'D17', 'D21', 'D22', 'D23', 'D25', 'D27', 'D28', 'D29', 'D30', 'D31', 'D32', 'D35', 'D36', 'D41', 'D45', 'D47', 'D49', 'D50', 'D51')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '44.246.98'
$ws.Range('E2').Value = '  +0.92%  '
$ws.Range('D3').Value = '2.244.22'
$ws.Range('E3').Value = '  +0.37%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').Value = '305.63'
$ws.Range('E5').Value = '  -2.83%  '
$ws.Range('D6').Value = '95.36'
$ws.Range('E6').Value = '  -2.89%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('E8').Value = '  +0.18%  '
$ws.Range('E9').Value = '  -1.36%  '
$ws.Range('D10').Value = '34.75'
$ws.Range('E10').Value = '  -2.69%  '
$ws.Range('E11').Value = '  -1.19%  '
$ws.Range('D12').Value = '7.20'
$ws.Range('E12').Value = '  -1.73%  '
$ws.Range('E13').Value = '  -0.15%  '
$ws.Range('D14').Value = '2.585.86'
$ws.Range('D15').Value = '2.328.37'
$ws.Range('E15').Value = '  +4.07%  '
$ws.Range('D16').Value = '0.831'
$ws.Range('E16').Value = '  -0.47%  '
$ws.Range('D17').Value = '13.54'
$ws.Range('E17').Value = '  -2.38%  '
$ws.Range('D18').Value = '44.012.06'
$ws.Range('E18').Value = '  +0.73%  '
$ws.Range('E19').Value = '  -0.50%  '
$ws.Range('E20').Value = '  +1.43%  '
$ws.Range('D21').Value = '12.03'
$ws.Range('E21').Value = '  -8.04%  '
$ws.Range('D22').Value = '65.48'
$ws.Range('E22').Value = '  -0.75%  '
$ws.Range('D23').Value = '237.67'
$ws.Range('E23').Value = '  +0.75%  '
$ws.Range('E24').Value = '  -1.37%  '
$ws.Range('D25').Value = '2.00'
$ws.Range('E25').Value = '  -1.23%  '
$ws.Range('E26').Value = '  -0.05%  '
$ws.Range('D27').Value = '9.91'
$ws.Range('E27').Value = '  -1.43%  '
$ws.Range('B28').Value = 'InjectiveProtocol'
$ws.Range('C28').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D28').Value = '37.92'
$ws.Range('E28').Value = '  +3.95%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').Value = '2.16'
$ws.Range('E29').Value = '  +1.13%  '
$ws.Range('D30').Value = '20.05'
$ws.Range('E30').Value = '  +0.36%  '
$ws.Range('D31').Value = '5.84'
$ws.Range('E31').Value = '  -1.91%  '
$ws.Range('D32').Value = '152.11'
$ws.Range('E32').Value = '  -2.37%  '
$ws.Range('E33').Value = '  -4.27%  '
$ws.Range('E34').Value = '  -0.86%  '
$ws.Range('D35').Value = '3.22'
$ws.Range('E35').Value = '  -3.07%  '
$ws.Range('D36').Value = '0.120'
$ws.Range('E36').Value = '  +2.08%  '
$ws.Range('E37').Value = '  -1.82%  '
$ws.Range('E38').Value = '  -7.43%  '
$ws.Range('E39').Value = '  +1.16%  '
$ws.Range('E40').Value = '  -3.90%  '
$ws.Range('D41').Value = '14.26'
$ws.Range('E41').Value = '  -8.47%  '
$ws.Range('E42').Value = '  -2.67%  '
$ws.Range('E43').Value = '  +0.22%  '
$ws.Range('D44').Value = '1.752.37'
$ws.Range('E44').Value = '  +3.17%  '
$ws.Range('D45').Value = '82.38'
$ws.Range('E45').Value = '  +0.11%  '
$ws.Range('E46').Value = '  -1.86%  '
$ws.Range('D47').Value = '99.85'
$ws.Range('E47').Value = '  -1.60%  '
$ws.Range('E48').Value = '  -4.02%  '
$ws.Range('D49').Value = '8.10'
$ws.Range('E49').Value = '  -0.35%  '
$ws.Range('D50').Value = '1.57'
$ws.Range('E50').Value = '  -1.97%  '
$ws.Range('D51').Value = '54.50'
$ws.Range('E51').Value = '  -2.92%  '
